# Apply CRC diagram updates for StockAdvisor design milestone.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 (ProgramRunner card: Responsibilities / Collaborators) ---
$ws.Range("A3").Value = "createAndShowGUI`naskStockTickerSymbol`naskInvestmentHorizon`nstartDataProcessing`nshowRecommendationResult"
$ws.Range("B3").Value = "UserInput`nUserOutput"
$ws.Range("D3").Value = "readStockTickerSymbol`nreadInvestmentHorizon"
$ws.Range("E3").Value = "UserInterface`nYahooFinanceDataReader"
$ws.Rows.Item(3).RowHeight = 90

# --- Row 6 (class-name headers) ---
$ws.Range("D6").Value = "ARIMAAnalysis"
$ws.Range("G6").Value = "ARMAAnalysis"

# --- Row 8 (Responsibilities / Collaborators under row 6 headers) ---
$ws.Range("A8").Value = "getStock`ngeneratePriceList`ngetHistoricalPrices"
$ws.Range("B8").Value = "UserInput`nDataProcessor"
$ws.Range("D8").Value = "runARIMA"
$ws.Range("E8").Value = "DataProcessor"
$ws.Range("G8").Value = "runARMA"
$ws.Range("H8").Value = "DataProcessor"
$ws.Rows.Item(8).RowHeight = 52.5

# --- Row 11 (class-name header) ---
$ws.Range("A11").Value = "DataProcessor"

# --- Row 13 (Responsibilities / Collaborators) ---
$ws.Range("A13").Value = "calculateAverage`ndataProcessing`ngetRecommendation"
$ws.Range("B13").Value = "UserInput`nYahooFinanceDataReader`nARIMAAnalysis`nARMAAnalysis`nUserOutput"
$ws.Range("D13").Value = "dataVisualization`ngiveRecommendation"
$ws.Range("E13").Value = "DataProcessor`nUserInterface"

# --- Row 16 (class-name header) ---
$ws.Range("A16").Value = "ProgramRunner"

# --- Row 18 ---
$ws.Range("A18").Value = "runStockAdvisor"

# --- Selection ---
$ws.Range("B13").Select()
